# Updated latest Guinea master data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean: clear all existing cell content/formatting on the sheet.
$ws.Cells.Clear()

# ---- Header row ----
$headers = @("id","name","descr","lang_code","is_active","cr_by","cr_dtimes","upd_by","upd_dtimes","is_deleted","del_dtimes")
for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---- Data rows ----
# Columns: A id | B name | C descr | D lang_code | E is_active | F cr_by | G cr_dtimes | H upd_by | I upd_dtimes | J is_deleted | K del_dtimes
$ids    = @(10001, 10003, 10005, 10007, 10009, 10011)
$names  = @("PrÃ©-Enregistrement", "Enregistrement", "Traitement", "Authentification", "Administration", "Portail RÃ©sident")
$descrs = @(
  "Portail Web pour les prÃ©-enregistrements",
  "Application pour les enregistrements",
  "Application pour les traitements post-enregistrements",
  "Application pour l'authentification des fournisseurs de services",
  "Portail Web pour la configuration des applications",
  "Portail Web pour les services dÃ©diÃ©s aux rÃ©sidents"
)

for ($i = 0; $i -lt $ids.Length; $i++) {
  $r = $i + 2
  $ws.Cells.Item($r, 1).Value = $ids[$i]
  $ws.Cells.Item($r, 2).Value = $names[$i]
  $ws.Cells.Item($r, 3).Value = $descrs[$i]
  $ws.Cells.Item($r, 4).Value = "fra"
  $ws.Cells.Item($r, 5).Value = $true
  $ws.Cells.Item($r, 6).Value = "superadmin"
  $ws.Cells.Item($r, 7).Value = 45079.576688067129
  $ws.Cells.Item($r, 7).NumberFormat = "mm:ss.0"
  $ws.Cells.Item($r, 8).Value = "NULL"
  $ws.Cells.Item($r, 9).Value = "NULL"
  $ws.Cells.Item($r, 10).Value = $false
  $ws.Cells.Item($r, 11).Value = "NULL"
}

# ---- Column widths ----
$ws.Columns.Item(4).ColumnWidth = 9.0

# ---- Selection ----
[void]$ws.Range("G14").Select()
